$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update kode_kecamatan (D2) and kode_desa (C2) sample values
$ws.Range("D2").Value = "KC1"
$ws.Range("C2").Value = "D01"

# Update the last active selection left over in the sheet view
$ws.Range("Q21").Select()
